# Refresh the cryptocurrency table (rows 2-51): latest prices and
# 1-hour volume changes, plus the newly-listed "Frax" row which shifts
# every following coin down by one row (pushing TrueUSD off the bottom).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price cells that are about to receive a
# numeric-looking string (e.g. "0.9946"), so Excel keeps storing them
# as text instead of silently converting them to numbers.
$ws.Range("D4:D9").NumberFormat = "@"
$ws.Range("D11:D16").NumberFormat = "@"
$ws.Range("D18:D20").NumberFormat = "@"
$ws.Range("D22:D25").NumberFormat = "@"
$ws.Range("D27:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '26.476.82'
$ws.Cells.Item(2, 5).Value = '  -0.48%  '

$ws.Cells.Item(3, 4).Value = '1.717.28'
$ws.Cells.Item(3, 5).Value = '  -1.49%  '

$ws.Cells.Item(4, 4).Value = '0.9946'
$ws.Cells.Item(4, 5).Value = '  -0.53%  '

$ws.Cells.Item(5, 4).Value = '240.25'
$ws.Cells.Item(5, 5).Value = '  -2.76%  '

$ws.Cells.Item(6, 4).Value = '0.9954'
$ws.Cells.Item(6, 5).Value = '  -0.49%  '

$ws.Cells.Item(7, 4).Value = '0.4889'
$ws.Cells.Item(7, 5).Value = '  -0.79%  '

$ws.Cells.Item(8, 4).Value = '0.2591'
$ws.Cells.Item(8, 5).Value = '  -3.31%  '

$ws.Cells.Item(9, 4).Value = '0.06185'
$ws.Cells.Item(9, 5).Value = '  -1.83%  '

$ws.Cells.Item(10, 4).Value = '1.723.71'
$ws.Cells.Item(10, 5).Value = '  -1.17%  '

$ws.Cells.Item(11, 4).Value = '0.06939'
$ws.Cells.Item(11, 5).Value = '  -1.68%  '

$ws.Cells.Item(12, 4).Value = '15.56'
$ws.Cells.Item(12, 5).Value = '  -1.22%  '

$ws.Cells.Item(13, 4).Value = '0.6021'
$ws.Cells.Item(13, 5).Value = '  -2.24%  '

$ws.Cells.Item(14, 4).Value = '4.462'
$ws.Cells.Item(14, 5).Value = '  -2.83%  '

$ws.Cells.Item(15, 4).Value = '76.59'
$ws.Cells.Item(15, 5).Value = '  -2.14%  '

$ws.Cells.Item(16, 4).Value = '0.9950'
$ws.Cells.Item(16, 5).Value = '  -0.51%  '

$ws.Cells.Item(17, 4).Value = '26.319.88'
$ws.Cells.Item(17, 5).Value = '  -1.09%  '

$ws.Cells.Item(18, 4).Value = '0.9947'
$ws.Cells.Item(18, 5).Value = '  -0.56%  '

$ws.Cells.Item(19, 4).Value = '0.000007098'
$ws.Cells.Item(19, 5).Value = '  -3.02%  '

$ws.Cells.Item(20, 4).Value = '11.26'
$ws.Cells.Item(20, 5).Value = '  -2.72%  '

$ws.Cells.Item(21, 4).Value = '1.935.36'

$ws.Cells.Item(22, 4).Value = '4.389'
$ws.Cells.Item(22, 5).Value = '  -4.33%  '

$ws.Cells.Item(23, 4).Value = '8.399'
$ws.Cells.Item(23, 5).Value = '  -3.98%  '

$ws.Cells.Item(24, 4).Value = '5.067'
$ws.Cells.Item(24, 5).Value = '  -4.02%  '

$ws.Cells.Item(25, 4).Value = '137.44'
$ws.Cells.Item(25, 5).Value = '  -1.53%  '

$ws.Cells.Item(26, 5).Value = '  -1.90%  '

$ws.Cells.Item(27, 4).Value = '1.395'
$ws.Cells.Item(27, 5).Value = '  -1.99%  '

$ws.Cells.Item(28, 4).Value = '1.736'
$ws.Cells.Item(28, 5).Value = '  -1.82%  '

$ws.Cells.Item(29, 4).Value = '105.39'
$ws.Cells.Item(29, 5).Value = '  -2.19%  '

$ws.Cells.Item(30, 4).Value = '3.899'
$ws.Cells.Item(30, 5).Value = '  -3.88%  '

$ws.Cells.Item(31, 4).Value = '0.07946'
$ws.Cells.Item(31, 5).Value = '  -1.39%  '

$ws.Cells.Item(32, 4).Value = '3.625'
$ws.Cells.Item(32, 5).Value = '  -3.13%  '

$ws.Cells.Item(33, 4).Value = '0.04461'
$ws.Cells.Item(33, 5).Value = '  -3.81%  '

$ws.Cells.Item(34, 2).Value = 'Frax'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(34, 4).Value = '0.9943'
$ws.Cells.Item(34, 5).Value = '  -0.51%  '

$ws.Cells.Item(35, 2).Value = 'HuobiToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(35, 4).Value = '2.588'
$ws.Cells.Item(35, 5).Value = '  -0.90%  '

$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(36, 4).Value = '0.9977'
$ws.Cells.Item(36, 5).Value = '  -2.24%  '

$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(37, 4).Value = '0.6155'
$ws.Cells.Item(37, 5).Value = '  -3.71%  '

$ws.Cells.Item(38, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(38, 4).Value = '0.9539'
$ws.Cells.Item(38, 5).Value = '  +5.88%  '

$ws.Cells.Item(39, 2).Value = 'RenderToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(39, 4).Value = '1.992'
$ws.Cells.Item(39, 5).Value = '  -3.64%  '

$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(40, 4).Value = '2.368'
$ws.Cells.Item(40, 5).Value = '  -2.46%  '

$ws.Cells.Item(41, 2).Value = 'PaxDollar'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(41, 4).Value = '0.9945'
$ws.Cells.Item(41, 5).Value = '  -0.91%  '

$ws.Cells.Item(42, 2).Value = 'VeChain'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(42, 4).Value = '0.01479'
$ws.Cells.Item(42, 5).Value = '  -1.80%  '

$ws.Cells.Item(43, 2).Value = 'Quant'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(43, 4).Value = '99.50'
$ws.Cells.Item(43, 5).Value = '  -2.41%  '

$ws.Cells.Item(44, 2).Value = 'FraxShare'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(44, 4).Value = '5.407'
$ws.Cells.Item(44, 5).Value = '  -0.55%  '

$ws.Cells.Item(45, 2).Value = 'TheSandbox'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(45, 4).Value = '0.3809'
$ws.Cells.Item(45, 5).Value = '  -3.17%  '

$ws.Cells.Item(46, 2).Value = 'Aptos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(46, 4).Value = '6.855'
$ws.Cells.Item(46, 5).Value = '  -0.80%  '

$ws.Cells.Item(47, 2).Value = 'Algorand'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(47, 4).Value = '0.1151'
$ws.Cells.Item(47, 5).Value = '  -2.92%  '

$ws.Cells.Item(48, 2).Value = 'Cronos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(48, 4).Value = '0.05351'
$ws.Cells.Item(48, 5).Value = '  -0.93%  '

$ws.Cells.Item(49, 2).Value = 'Elrond'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Cells.Item(49, 4).Value = '30.36'
$ws.Cells.Item(49, 5).Value = '  -0.78%  '

$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '7.711'
$ws.Cells.Item(50, 5).Value = '  -1.94%  '

$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).Value = '51.15'
$ws.Cells.Item(51, 5).Value = '  -1.19%  '
